$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 13.32347289536497
$ws.Cells.Item(2, 3).Value = 8.621370384004305
$ws.Cells.Item(2, 4).Value = 5.976175885528134
$ws.Cells.Item(2, 5).Value = 11.13812739258253
$ws.Cells.Item(2, 7).Value = 3.620795274893379
$ws.Cells.Item(2, 13).Value = 14.94262488527407
$ws.Cells.Item(2, 15).Value = 21.92193066668024

$ws.Cells.Item(3, 2).Value = 12.72803815671297
$ws.Cells.Item(3, 3).Value = 8.146525002225319
$ws.Cells.Item(3, 4).Value = 5.856247262053571
$ws.Cells.Item(3, 5).Value = 11.04363809332773
$ws.Cells.Item(3, 7).Value = 3.623560066773948
$ws.Cells.Item(3, 13).Value = 14.65559845052419
$ws.Cells.Item(3, 15).Value = 21.93792121216747

$ws.Cells.Item(4, 2).Value = 12.34957036039175
$ws.Cells.Item(4, 3).Value = 7.838918388945147
$ws.Cells.Item(4, 4).Value = 5.783220541135611
$ws.Cells.Item(4, 5).Value = 10.98976343518347
$ws.Cells.Item(4, 7).Value = 3.625345876694209
$ws.Cells.Item(4, 13).Value = 14.48007159145226
$ws.Cells.Item(4, 15).Value = 21.9555624881631

$ws.Cells.Item(5, 2).Value = 12.19232259089783
$ws.Cells.Item(5, 3).Value = 7.709587761525345
$ws.Cells.Item(5, 4).Value = 5.753660658480785
$ws.Cells.Item(5, 5).Value = 10.96887064603442
$ws.Cells.Item(5, 7).Value = 3.626095869763923
$ws.Cells.Item(5, 13).Value = 14.40882133633727
$ws.Cells.Item(5, 15).Value = 21.96470938720958

$ws.Cells.Item(6, 2).Value = 12.16603630453172
$ws.Cells.Item(6, 3).Value = 7.687874162496983
$ws.Cells.Item(6, 4).Value = 5.748765630639035
$ws.Cells.Item(6, 5).Value = 10.96546606357543
$ws.Cells.Item(6, 7).Value = 3.626221752352669
$ws.Cells.Item(6, 13).Value = 14.39701004852023
$ws.Cells.Item(6, 15).Value = 21.96634617866205

$ws.Cells.Item(7, 2).Value = 12.34746157744608
$ws.Cells.Item(7, 3).Value = 7.837190213854359
$ws.Cells.Item(7, 4).Value = 5.78282101937467
$ws.Cells.Item(7, 5).Value = 10.98947734499465
$ws.Cells.Item(7, 7).Value = 3.625355901126004
$ws.Cells.Item(7, 13).Value = 14.47910942528333
$ws.Cells.Item(7, 15).Value = 21.95567793222105

$ws.Cells.Item(8, 2).Value = 13.12095879765254
$ws.Cells.Item(8, 3).Value = 8.461005274924783
$ws.Cells.Item(8, 4).Value = 5.934725191188106
$ws.Cells.Item(8, 5).Value = 11.10470056230735
$ws.Cells.Item(8, 7).Value = 3.621730313000457
$ws.Cells.Item(8, 13).Value = 14.84357080776323
$ws.Cells.Item(8, 15).Value = 21.92581544034024

$ws.Cells.Item(9, 2).Value = 14.52760099653715
$ws.Cells.Item(9, 3).Value = 9.55504584469228
$ws.Cells.Item(9, 4).Value = 6.235471198971563
$ws.Cells.Item(9, 5).Value = 11.3625172448538
$ws.Cells.Item(9, 7).Value = 3.6153169301387
$ws.Cells.Item(9, 13).Value = 15.55958298685957
$ws.Cells.Item(9, 15).Value = 21.92966982485933

$ws.Cells.Item(10, 2).Value = 15.48487451475165
$ws.Cells.Item(10, 3).Value = 10.27816656473774
$ws.Cells.Item(10, 4).Value = 6.455603392443648
$ws.Cells.Item(10, 5).Value = 11.56986023531366
$ws.Cells.Item(10, 7).Value = 3.611024545042937
$ws.Cells.Item(10, 13).Value = 16.08083129027146
$ws.Cells.Item(10, 15).Value = 21.97091519447151

$ws.Cells.Item(11, 2).Value = 15.90236364527738
$ws.Cells.Item(11, 3).Value = 10.58937515305615
$ws.Cells.Item(11, 4).Value = 6.555054490055859
$ws.Cells.Item(11, 5).Value = 11.66771923156803
$ws.Cells.Item(11, 7).Value = 3.609161855976227
$ws.Cells.Item(11, 13).Value = 16.3157996192485
$ws.Cells.Item(11, 15).Value = 21.99806272039681

$ws.Cells.Item(12, 2).Value = 16.05777397573809
$ws.Cells.Item(12, 3).Value = 10.7046572841198
$ws.Cells.Item(12, 4).Value = 6.592574978141215
$ws.Cells.Item(12, 5).Value = 11.70525234249934
$ws.Cells.Item(12, 7).Value = 3.608469355560975
$ws.Cells.Item(12, 13).Value = 16.40437940669069
$ws.Cells.Item(12, 15).Value = 22.00954914098037

$ws.Cells.Item(13, 2).Value = 16.02442433733088
$ws.Cells.Item(13, 3).Value = 10.67994351096658
$ws.Cells.Item(13, 4).Value = 6.584501094635393
$ws.Cells.Item(13, 5).Value = 11.69714827534241
$ws.Cells.Item(13, 7).Value = 3.60861792715246
$ws.Cells.Item(13, 13).Value = 16.38532120557958
$ws.Cells.Item(13, 15).Value = 22.00702169767182

$ws.Cells.Item(14, 2).Value = 15.91520363938873
$ws.Cells.Item(14, 3).Value = 10.59891097237741
$ws.Cells.Item(14, 4).Value = 6.558144342803073
$ws.Cells.Item(14, 5).Value = 11.67079775374166
$ws.Cells.Item(14, 7).Value = 3.609104626246215
$ws.Cells.Item(14, 13).Value = 16.32309558997594
$ws.Cells.Item(14, 15).Value = 21.99898354366103

$ws.Cells.Item(15, 2).Value = 15.84795066125365
$ws.Cells.Item(15, 3).Value = 10.54894170913388
$ws.Cells.Item(15, 4).Value = 6.54198074179835
$ws.Cells.Item(15, 5).Value = 11.6547183037255
$ws.Cells.Item(15, 7).Value = 3.609404416006311
$ws.Cells.Item(15, 13).Value = 16.28492627932955
$ws.Cells.Item(15, 15).Value = 21.99421701527522

$ws.Cells.Item(16, 2).Value = 15.45722004924359
$ws.Cells.Item(16, 3).Value = 10.25746998240736
$ws.Cells.Item(16, 4).Value = 6.449086666661474
$ws.Cells.Item(16, 5).Value = 11.56353333857317
$ws.Cells.Item(16, 7).Value = 3.611148079650536
$ws.Cells.Item(16, 13).Value = 16.06542462299382
$ws.Cells.Item(16, 15).Value = 21.96930984855968

$ws.Cells.Item(17, 2).Value = 15.21283931063702
$ws.Cells.Item(17, 3).Value = 10.07410540225128
$ws.Cells.Item(17, 4).Value = 6.39189337205288
$ws.Cells.Item(17, 5).Value = 11.50847750902329
$ws.Cells.Item(17, 7).Value = 3.61224074376655
$ws.Cells.Item(17, 13).Value = 15.9301523449671
$ws.Cells.Item(17, 15).Value = 21.95617859360484

$ws.Cells.Item(18, 2).Value = 15.07059042874944
$ws.Cells.Item(18, 3).Value = 9.966969452764522
$ws.Cells.Item(18, 4).Value = 6.35893481907468
$ws.Cells.Item(18, 5).Value = 11.47714586841105
$ws.Cells.Item(18, 7).Value = 3.612877685065347
$ws.Cells.Item(18, 13).Value = 15.85214992517561
$ws.Cells.Item(18, 15).Value = 21.94941514238973

$ws.Cells.Item(19, 2).Value = 15.02214082908068
$ws.Cells.Item(19, 3).Value = 9.930408841904825
$ws.Cells.Item(19, 4).Value = 6.347766159179892
$ws.Cells.Item(19, 5).Value = 11.46659602943596
$ws.Cells.Item(19, 7).Value = 3.613094799389451
$ws.Cells.Item(19, 13).Value = 15.82570851755928
$ws.Cells.Item(19, 15).Value = 21.94726067194292

$ws.Cells.Item(20, 2).Value = 15.23902947558013
$ws.Cells.Item(20, 3).Value = 10.0937977804074
$ws.Cells.Item(20, 4).Value = 6.397988463327768
$ws.Cells.Item(20, 5).Value = 11.51430385628473
$ws.Cells.Item(20, 7).Value = 3.612123551747262
$ws.Cells.Item(20, 13).Value = 15.94457339405593
$ws.Cells.Item(20, 15).Value = 21.95749473907356

$ws.Cells.Item(21, 2).Value = 15.94735794643261
$ws.Cells.Item(21, 3).Value = 10.62278192083632
$ws.Cells.Item(21, 4).Value = 6.56589005710794
$ws.Cells.Item(21, 5).Value = 11.67852488698868
$ws.Cells.Item(21, 7).Value = 3.60896132250499
$ws.Cells.Item(21, 13).Value = 16.34138423382139
$ws.Cells.Item(21, 15).Value = 22.00131181204341

$ws.Cells.Item(22, 2).Value = 16.39461386861253
$ws.Cells.Item(22, 3).Value = 10.95354451709398
$ws.Cells.Item(22, 4).Value = 6.674793250688162
$ws.Cells.Item(22, 5).Value = 11.78861166081492
$ws.Cells.Item(22, 7).Value = 3.606969544116318
$ws.Cells.Item(22, 13).Value = 16.59836785485784
$ws.Cells.Item(22, 15).Value = 22.03697867591324

$ws.Cells.Item(23, 2).Value = 16.15736792717033
$ws.Cells.Item(23, 3).Value = 10.77838267648813
$ws.Cells.Item(23, 4).Value = 6.616758332200304
$ws.Cells.Item(23, 5).Value = 11.7296149005998
$ws.Cells.Item(23, 7).Value = 3.608025762622559
$ws.Cells.Item(23, 13).Value = 16.46145432133601
$ws.Cells.Item(23, 15).Value = 22.01729965790006

$ws.Cells.Item(24, 2).Value = 15.22719434880961
$ws.Cells.Item(24, 3).Value = 10.08490020395358
$ws.Cells.Item(24, 4).Value = 6.395233110472444
$ws.Cells.Item(24, 5).Value = 11.51166876331853
$ws.Cells.Item(24, 7).Value = 3.612176506983225
$ws.Cells.Item(24, 13).Value = 15.93805435543498
$ws.Cells.Item(24, 15).Value = 21.95689726166601

$ws.Cells.Item(25, 2).Value = 14.1598850116028
$ws.Cells.Item(25, 3).Value = 9.273142199329779
$ws.Cells.Item(25, 4).Value = 6.154076012609985
$ws.Cells.Item(25, 5).Value = 11.28950660741265
$ws.Cells.Item(25, 7).Value = 3.616977886516604
$ws.Cells.Item(25, 13).Value = 15.36633387195918
$ws.Cells.Item(25, 15).Value = 21.92190105532656

